# edit.ps1 - apply update-log diff changes via Word COM interop
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "-present development." -> "-5/19/2022 (deployed to shiny
# server 5/19/22).", keeping the surrounding "11/28/2021" run and the
# trailing space run separate (5 runs total in the rebuilt paragraph).
# ---------------------------------------------------------------------
$target1 = "11/28/2021-present development. "
$full = $d.Content.Text
$idx1 = $full.IndexOf($target1)
if ($idx1 -lt 0) {
    throw "Could not locate target text for change 1"
}
$rng1 = $d.Range($idx1, $idx1 + $target1.Length)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>11/28/2021</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>5/19/2022 (deployed to shiny server 5/19/22)</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2: append three new bullet paragraphs after the "Modified user
# manual..." paragraph (6/8/2022 note + sub-note, and 9/1/2022 note).
# ---------------------------------------------------------------------
$marker2 = "Modified user manual to detail percentile tabs and account for several of the above changes."
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf($marker2)
if ($idx2 -lt 0) {
    throw "Could not locate target text for change 2"
}
$endIdx2 = $idx2 + $marker2.Length
$rng2 = $d.Range($endIdx2, $endIdx2)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:textAlignment w:val="baseline"/></w:pPr><w:r><w:t>6/8/2022 Added checkbox and code that allows ALK to extrapolate to fish smaller than those aged.  Saw that ODWC frequently had a floor on fish sizes that are aged (e.g., crappie less than 100mm are always age-0 and are never aged).  This produced inaccurate length at age 0 as only fish 100mm and larger were considered in the calculation given that was the smallest size aged.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:textAlignment w:val="baseline"/></w:pPr><w:r><w:t>Default is to allow this.  Only possible problem would be if only age-1+ are aged…then the check box might try to parse age-1 vs age-0 incorrectly (or might assume all smaller fish are age-1 and no age-0 exist)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:textAlignment w:val="baseline"/></w:pPr><w:r><w:t>9/1/2022</w:t></w:r><w:r><w:t xml:space="preserve"> – Modified to account for Verified.TL and Verified.Wr files in user-uploaded files</w:t></w:r><w:r><w:t xml:space="preserve"> rather than common Verified.TL.Wr column (was afraid users would fix TL issues but miss Wr issues that are left related to weight errors)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2.InsertXML($xml2)
